# Auto-generated update of leve-profit market data columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.
# For each affected row we only touch the cells whose value actually
# changed; cells that disappear in the target state are cleared so the
# underlying <c> node is removed (matching a refreshed/blank market quote).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 98
$ws.Range("H98").Value = 1562
$ws.Range("I98").Value = 1321.3334
$ws.Range("J98").Value = 3006
$ws.Range("K98").Value = 1321.3334
$ws.Range("L98").Value = 3006
$ws.Range("M98").Value = 176.6666
$ws.Range("N98").Value = -6002
# Row 122
$ws.Range("H122").Value = 1562
$ws.Range("I122").Value = 1321.3334
$ws.Range("J122").Value = 3006
$ws.Range("K122").Value = 3964.0002
$ws.Range("L122").Value = 9018
$ws.Range("M122").Value = -1514.0002
$ws.Range("N122").Value = -13918
# Row 129
$ws.Range("H129").Value = 1024.7142
$ws.Range("I129").Value = 601.3043
$ws.Range("J129").Value = 2972.4
$ws.Range("K129").Value = 1803.9129
$ws.Range("L129").Value = 8917.200000000001
$ws.Range("M129").Value = 3196.0871
$ws.Range("N129").Value = -18917.2
# Row 132
$ws.Range("H132").Value = 2058.8262
$ws.Range("I132").Value = 2088.7727
$ws.Range("K132").Value = 6266.3181
$ws.Range("M132").Value = -3736.3181
# Row 135
$ws.Range("H135").Value = 966.3333
$ws.Range("I135").Value = 966.3333
$ws.Range("K135").Value = 8696.9997
$ws.Range("M135").Value = -6161.9997
# Row 138
$ws.Range("H138").Value = 3721.1924
$ws.Range("I138").Value = 2146.75
$ws.Range("J138").Value = 4007.4546
$ws.Range("K138").Value = 6440.25
$ws.Range("L138").Value = 12022.3638
$ws.Range("M138").Value = -1300.25
$ws.Range("N138").Value = -22302.3638

$ws = $wb.Worksheets.Item("ARM")

# Row 74
$ws.Range("H74").Value = 3412.2727
$ws.Range("I74").Value = 3506.5
$ws.Range("K74").Value = 3506.5
$ws.Range("M74").Value = -2632.5
# Row 77
$ws.Range("H77").Value = 3412.2727
$ws.Range("I77").Value = 3506.5
$ws.Range("K77").Value = 17532.5
$ws.Range("M77").Value = -13164.5
# Row 88
$ws.Range("H88").Value = 2918.7334
$ws.Range("J88").Value = 3370.9167
$ws.Range("L88").Value = 3370.9167
$ws.Range("N88").Value = -4182.9167
# Row 91
$ws.Range("H91").Value = 2918.7334
$ws.Range("J91").Value = 3370.9167
$ws.Range("L91").Value = 3370.9167
$ws.Range("N91").Value = -6178.9167
# Row 122
$ws.Range("H122").Value = 1500
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()
# Row 132
$ws.Range("H132").Value = 1368.1111
$ws.Range("I132").Value = 1247.2727
$ws.Range("K132").Value = 3741.8181
$ws.Range("M132").Value = -1211.8181

$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 3018.6667
$ws.Range("I86").Value = 3018.6667
$ws.Range("K86").Value = 3018.6667
$ws.Range("M86").Value = -1895.6667
# Row 89
$ws.Range("H89").Value = 3018.6667
$ws.Range("I89").Value = 3018.6667
$ws.Range("K89").Value = 15093.3335
$ws.Range("M89").Value = -9477.333500000001
# Row 134
$ws.Range("H134").Value = 4245.357
$ws.Range("I134").Value = 5792.8887
$ws.Range("J134").Value = 1459.8
$ws.Range("K134").Value = 17378.6661
$ws.Range("L134").Value = 4379.4
$ws.Range("M134").Value = -14843.6661
$ws.Range("N134").Value = -9449.4
# Row 135
$ws.Range("H135").Value = 49997.6
$ws.Range("J135").Value = 49997.6
$ws.Range("L135").Value = 49997.6
$ws.Range("N135").Value = -60137.6

$ws = $wb.Worksheets.Item("CRP")

# Row 86
$ws.Range("H86").Value = 10102
$ws.Range("I86").Value = 10165.2
$ws.Range("K86").Value = 10165.2
$ws.Range("M86").Value = -9042.200000000001
# Row 89
$ws.Range("H89").Value = 10102
$ws.Range("I89").Value = 10165.2
$ws.Range("K89").Value = 50826
$ws.Range("M89").Value = -45210
# Row 105
$ws.Range("H105").Value = 1294.625
$ws.Range("I105").Value = 1265.2858
$ws.Range("K105").Value = 1265.2858
$ws.Range("M105").Value = 481.7141999999999

$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 2067.5
$ws.Range("I4").Value = 2281.25
$ws.Range("K4").Value = 6843.75
$ws.Range("M4").Value = -6731.75
# Row 13
$ws.Range("H13").Value = 289
$ws.Range("I13").Value = 649.6667
$ws.Range("J13").Value = 18.5
$ws.Range("K13").Value = 1949.0001
$ws.Range("L13").Value = 55.5
$ws.Range("M13").Value = -1781.0001
$ws.Range("N13").Value = -391.5
# Row 122
$ws.Range("H122").Value = 471.875
$ws.Range("J122").Value = 457.2
$ws.Range("L122").Value = 4114.8
$ws.Range("N122").Value = -9014.799999999999
# Row 129
$ws.Range("H129").Value = 2257.6667
$ws.Range("I129").Value = 388
$ws.Range("K129").Value = 1164
$ws.Range("M129").Value = 3836

$ws = $wb.Worksheets.Item("GSM")

# Row 113
$ws.Range("H113").Value = 3374.8
$ws.Range("I113").Value = 3593.5
$ws.Range("K113").Value = 3593.5
$ws.Range("M113").Value = -1423.5

$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Range("H7").Value = 7599
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
# Row 16
$ws.Range("H16").Value = 3397.6
$ws.Range("I16").Value = 3497
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 3497
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -3327
$ws.Range("N16").Value = -3340
# Row 68
$ws.Range("H68").Value = 58500
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 58500
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 58500
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -59998
# Row 71
$ws.Range("H71").Value = 58500
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 58500
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 292500
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -299988
# Row 82
$ws.Range("H82").Value = 2253.125
$ws.Range("I82").Value = 1421.75
$ws.Range("J82").Value = 3084.5
$ws.Range("K82").Value = 1421.75
$ws.Range("L82").Value = 3084.5
$ws.Range("M82").Value = -1060.75
$ws.Range("N82").Value = -3806.5
# Row 85
$ws.Range("H85").Value = 2253.125
$ws.Range("I85").Value = 1421.75
$ws.Range("J85").Value = 3084.5
$ws.Range("K85").Value = 1421.75
$ws.Range("L85").Value = 3084.5
$ws.Range("M85").Value = -173.75
$ws.Range("N85").Value = -5580.5
# Row 126
$ws.Range("H126").Value = 7599
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
# Row 132
$ws.Range("H132").Value = 6449.533
$ws.Range("I132").Value = 4881.4
$ws.Range("J132").Value = 7233.6
$ws.Range("K132").Value = 14644.2
$ws.Range("L132").Value = 21700.8
$ws.Range("M132").Value = -12114.2
$ws.Range("N132").Value = -26760.8
# Row 136
$ws.Range("H136").Value = 3408.75
$ws.Range("I136").Value = 3213.5
$ws.Range("J136").Value = 3994.5
$ws.Range("K136").Value = 9640.5
$ws.Range("L136").Value = 11983.5
$ws.Range("M136").Value = -7090.5
$ws.Range("N136").Value = -17083.5

$ws = $wb.Worksheets.Item("WVR")

# Row 59
$ws.Range("H59").Value = 1500
$ws.Range("J59").Value = 1500
$ws.Range("L59").Value = 1500
$ws.Range("N59").Value = -2976
# Row 81
$ws.Range("H81").Value = 4397.125
$ws.Range("J81").Value = 4795.25
$ws.Range("L81").Value = 9590.5
$ws.Range("N81").Value = -11712.5
# Row 84
$ws.Range("H84").Value = 4397.125
$ws.Range("J84").Value = 4795.25
$ws.Range("L84").Value = 47952.5
$ws.Range("N84").Value = -58560.5
# Row 122
$ws.Range("H122").Value = 3441
$ws.Range("I122").Value = 3500
$ws.Range("J122").Value = 3426.25
$ws.Range("K122").Value = 10500
$ws.Range("L122").Value = 10278.75
$ws.Range("M122").Value = -8050
$ws.Range("N122").Value = -15178.75
# Row 127
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
# Row 136
$ws.Range("H136").Value = 6134.8423
$ws.Range("I136").Value = 6151.8667
$ws.Range("J136").Value = 6071
$ws.Range("K136").Value = 18455.6001
$ws.Range("L136").Value = 18213
$ws.Range("M136").Value = -15905.6001
$ws.Range("N136").Value = -23313
